# Scheduled price-data refresh for the Anima_Profits crafting-leve profit sheets.
# For each (sheet, row) the currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) are refreshed from the latest market-board snapshot.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1126.3846
$ws.Range("I129").Value = 651.5454999999999
$ws.Range("J129").Value = 1312.9286
$ws.Range("K129").Value = 1954.6365
$ws.Range("L129").Value = 3938.7858
$ws.Range("M129").Value = 3045.3635
$ws.Range("N129").Value = -13938.7858
$ws.Range("H138").Value = 2606.1
$ws.Range("I138").Value = 2832.7144
$ws.Range("J138").Value = 2517.9722
$ws.Range("K138").Value = 8498.143199999999
$ws.Range("L138").Value = 7553.9166
$ws.Range("M138").Value = -3358.143199999999
$ws.Range("N138").Value = -17833.9166
$ws.Range("H141").Value = 5170.6523
$ws.Range("I141").Value = 2127.625
$ws.Range("J141").Value = 12126.143
$ws.Range("K141").Value = 6382.875
$ws.Range("L141").Value = 36378.429
$ws.Range("M141").Value = -1202.875
$ws.Range("N141").Value = -46738.429

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 319236.22
$ws.Range("I32").Value = 369475.72
$ws.Range("J32").Value = 51292.266
$ws.Range("K32").Value = 369475.72
$ws.Range("L32").Value = 51292.266
$ws.Range("M32").Value = -369188.72
$ws.Range("N32").Value = -51866.266
$ws.Range("H43").Value = 13999.5
$ws.Range("J43").Value = 13000
$ws.Range("L43").Value = 13000
$ws.Range("N43").Value = -13626
$ws.Range("H52").Value = 98780
$ws.Range("J52").Value = 98780
$ws.Range("L52").Value = 98780
$ws.Range("N52").Value = -99416
$ws.Range("H61").Value = 7938792.5
$ws.Range("I61").Value = 20834734
$ws.Range("J61").Value = 2829
$ws.Range("K61").Value = 20834734
$ws.Range("L61").Value = 2829
$ws.Range("M61").Value = -20834522
$ws.Range("N61").Value = -3253
$ws.Range("H74").Value = 841.4074000000001
$ws.Range("I74").Value = 566.6667
$ws.Range("J74").Value = 1061.2
$ws.Range("K74").Value = 566.6667
$ws.Range("L74").Value = 1061.2
$ws.Range("M74").Value = 307.3333
$ws.Range("N74").Value = -2809.2
$ws.Range("H77").Value = 841.4074000000001
$ws.Range("I77").Value = 566.6667
$ws.Range("J77").Value = 1061.2
$ws.Range("K77").Value = 2833.3335
$ws.Range("L77").Value = 5306
$ws.Range("M77").Value = 1534.6665
$ws.Range("N77").Value = -14042
$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -1150
$ws.Range("H132").Value = 5597.4062
$ws.Range("I132").Value = 6663.421
$ws.Range("J132").Value = 4039.3845
$ws.Range("K132").Value = 19990.263
$ws.Range("L132").Value = 12118.1535
$ws.Range("M132").Value = -17460.263
$ws.Range("N132").Value = -17178.1535
$ws.Range("H136").Value = 7938792.5
$ws.Range("I136").Value = 20834734
$ws.Range("J136").Value = 2829
$ws.Range("K136").Value = 62504202
$ws.Range("L136").Value = 8487
$ws.Range("M136").Value = -62501652
$ws.Range("N136").Value = -13587

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 887.4545000000001
$ws.Range("I94").Value = 770.5
$ws.Range("J94").Value = 1199.3334
$ws.Range("K94").Value = 770.5
$ws.Range("L94").Value = 1199.3334
$ws.Range("M94").Value = -319.5
$ws.Range("N94").Value = -2101.3334
$ws.Range("H99").Value = 1659
$ws.Range("I99").Value = 1659
$ws.Range("K99").Value = 1659
$ws.Range("M99").Value = -161
$ws.Range("H134").Value = 2594.5676
$ws.Range("I134").Value = 2287.5518
$ws.Range("J134").Value = 3707.5
$ws.Range("K134").Value = 6862.655400000001
$ws.Range("L134").Value = 11122.5
$ws.Range("M134").Value = -4327.655400000001
$ws.Range("N134").Value = -16192.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 40001.5
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 70003
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 70003
$ws.Range("M3").Value = -9887
$ws.Range("N3").Value = -70229
$ws.Range("H31").Value = 5040.3335
$ws.Range("I31").Value = 1209.7916
$ws.Range("J31").Value = 12701.417
$ws.Range("K31").Value = 1209.7916
$ws.Range("L31").Value = 12701.417
$ws.Range("M31").Value = -914.7916
$ws.Range("N31").Value = -13291.417
$ws.Range("H34").Value = 5040.3335
$ws.Range("I34").Value = 1209.7916
$ws.Range("J34").Value = 12701.417
$ws.Range("K34").Value = 1209.7916
$ws.Range("L34").Value = 12701.417
$ws.Range("M34").Value = -1007.7916
$ws.Range("N34").Value = -13105.417
$ws.Range("H44").Value = 18035.5
$ws.Range("J44").Value = 18035.5
$ws.Range("L44").Value = 18035.5
$ws.Range("N44").Value = -18919.5
$ws.Range("H58").Value = 2463.3333
$ws.Range("I58").Value = 2426.625
$ws.Range("J58").Value = 2757
$ws.Range("K58").Value = 2426.625
$ws.Range("L58").Value = 2757
$ws.Range("M58").Value = -2223.625
$ws.Range("N58").Value = -3163
$ws.Range("H136").Value = 2463.3333
$ws.Range("I136").Value = 2426.625
$ws.Range("J136").Value = 2757
$ws.Range("K136").Value = 7279.875
$ws.Range("L136").Value = 8271
$ws.Range("M136").Value = -4729.875
$ws.Range("N136").Value = -13371

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 124.789474
$ws.Range("I2").Value = 52.333332
$ws.Range("J2").Value = 190
$ws.Range("K2").Value = 313.999992
$ws.Range("L2").Value = 1140
$ws.Range("M2").Value = -200.999992
$ws.Range("N2").Value = -1366
$ws.Range("H26").Value = 373.64
$ws.Range("I26").Value = 48.714287
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 146.142861
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = 141.857139
$ws.Range("N26").Value = -2076
$ws.Range("H40").Value = 411
$ws.Range("I40").Value = 411
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1644
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1575
$ws.Range("H68").Value = 1382.1305
$ws.Range("I68").Value = 817.7143
$ws.Range("J68").Value = 1549.0704
$ws.Range("K68").Value = 2453.1429
$ws.Range("L68").Value = 4647.2112
$ws.Range("M68").Value = -1642.1429
$ws.Range("N68").Value = -6269.2112
$ws.Range("H71").Value = 1382.1305
$ws.Range("I71").Value = 817.7143
$ws.Range("J71").Value = 1549.0704
$ws.Range("K71").Value = 7359.428699999999
$ws.Range("L71").Value = 13941.6336
$ws.Range("M71").Value = -3303.428699999999
$ws.Range("N71").Value = -22053.6336
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("N101").Value = 0
$ws.Range("H131").Value = 1239.375
$ws.Range("I131").Value = 1030
$ws.Range("J131").Value = 1253.3334
$ws.Range("K131").Value = 3090
$ws.Range("L131").Value = 3760.0002
$ws.Range("M131").Value = 1950
$ws.Range("N131").Value = -13840.0002

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1333.8182
$ws.Range("I97").Value = 1465.875
$ws.Range("J97").Value = 981.6667
$ws.Range("K97").Value = 1465.875
$ws.Range("L97").Value = 981.6667
$ws.Range("M97").Value = -969.875
$ws.Range("N97").Value = -1973.6667
$ws.Range("H122").Value = 4109.5713
$ws.Range("I122").Value = 3326.625
$ws.Range("J122").Value = 4591.385
$ws.Range("K122").Value = 9979.875
$ws.Range("L122").Value = 13774.155
$ws.Range("M122").Value = -7529.875
$ws.Range("N122").Value = -18674.155
$ws.Range("H132").Value = 2230.8125
$ws.Range("I132").Value = 1866.2222
$ws.Range("K132").Value = 5598.6666
$ws.Range("M132").Value = -3068.6666

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2234.111
$ws.Range("I93").Value = 1300
$ws.Range("J93").Value = 2701.1667
$ws.Range("K93").Value = 1300
$ws.Range("L93").Value = 2701.1667
$ws.Range("M93").Value = -52
$ws.Range("N93").Value = -5197.1667
$ws.Range("H132").Value = 4759.1816
$ws.Range("I132").Value = 4985.6665
$ws.Range("J132").Value = 4487.4
$ws.Range("K132").Value = 14956.9995
$ws.Range("L132").Value = 13462.2
$ws.Range("M132").Value = -12426.9995
$ws.Range("N132").Value = -18522.2
$ws.Range("H136").Value = 6174128.5
$ws.Range("I136").Value = 1349
$ws.Range("J136").Value = 15152717
$ws.Range("K136").Value = 4047
$ws.Range("L136").Value = 45458151
$ws.Range("M136").Value = -1497
$ws.Range("N136").Value = -45463251

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 26666.666
$ws.Range("J48").Value = 26666.666
$ws.Range("L48").Value = 26666.666
$ws.Range("N48").Value = -27804.666
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0
$ws.Range("H100").Value = 1891.6666
$ws.Range("I100").Value = 850
$ws.Range("J100").Value = 2933.3333
$ws.Range("K100").Value = 1700
$ws.Range("L100").Value = 5866.6666
$ws.Range("M100").Value = -1159
$ws.Range("N100").Value = -6948.6666
$ws.Range("H122").Value = 1950.3334
$ws.Range("I122").Value = 1450.5
$ws.Range("K122").Value = 4351.5
$ws.Range("M122").Value = -1901.5
$ws.Range("H136").Value = 3092.077
$ws.Range("I136").Value = 2999.8262
$ws.Range("J136").Value = 3224.6875
$ws.Range("K136").Value = 8999.4786
$ws.Range("L136").Value = 9674.0625
$ws.Range("M136").Value = -6449.4786
$ws.Range("N136").Value = -14774.0625

